$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at AF, shifting the existing "area" column (AF) to AG.
$ws.Columns("AF").Insert()

# New header for the inserted column
$ws.Range("AF1").Value = "interval"

# Populate the new "interval" column for data rows 2-180 with "inside"/"outside"
# based on whether each sample falls inside or outside the plotted interval.
$ws.Range("AF2:AF18").Value = "inside"
$ws.Range("AF19").Value = "outside"
$ws.Range("AF20:AF29").Value = "inside"
$ws.Range("AF30").Value = "outside"
$ws.Range("AF31:AF59").Value = "inside"
$ws.Range("AF60").Value = "outside"
$ws.Range("AF61:AF65").Value = "inside"
$ws.Range("AF66").Value = "outside"
$ws.Range("AF67:AF159").Value = "inside"
$ws.Range("AF160").Value = "outside"
$ws.Range("AF161:AF163").Value = "inside"
$ws.Range("AF164").Value = "outside"
$ws.Range("AF165:AF168").Value = "inside"
$ws.Range("AF169").Value = "outside"
$ws.Range("AF170:AF180").Value = "inside"
